# Generate Report for Handback
#
# Refreshes the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps for the first data
# row (the "ab276e7c..." file) on both the "zh-cn" and "de-de" sheets, as
# would happen when the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-27 02:35:19"
$zhcn.Range("G2").Value = "2016-01-27 02:36:09"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-27 02:35:31"
$dede.Range("G2").Value = "2016-01-27 02:36:31"
